# Apply updated "dSF" (column F) values as per repulled data / recalculated mean.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = 4
    4  = -3
    5  = 3
    10 = 5
    11 = -1
    14 = -6
    17 = 2
    18 = -2
    19 = 3
    20 = -4
    21 = -2
    22 = 5
    23 = 4
    24 = 2
    25 = 1
    26 = 2
    27 = -1
    28 = 1
    33 = -4
    34 = -3
    36 = 1
    38 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
